$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Row 20 ("ELITE LIQUOR INC") now has a Last Invoice Date.
#    Copy the date-format style from an existing dated cell (D2) so the
#    new value picks up style index 12 (numFmt 165, left/top aligned)
#    instead of creating a brand-new style entry.
# -----------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 45954

# -----------------------------------------------------------------
# 2) A brand-new prospect ("ROSALIA LLC") was inserted right before the
#    "HOLY FAMILY MARONITE CHURCH" row, pushing that row (and everything
#    that used to be below the insertion point) down by one.
# -----------------------------------------------------------------
$ws.Rows(23).Insert()
$ws.Rows(23).RowHeight = 13.05

$ws.Range("A23").Value = "ROSALIA LLC"
$ws.Range("B23").Value = "Cina, Jonathan D"
$ws.Range("C23").Value = "023"
$ws.Range("E23").Value = "0008350"

# -----------------------------------------------------------------
# 3) A new row was appended at the bottom of the table for
#    "SCHMITT MUSIC CTR". Copy formatting from the row directly above
#    (the Holy Family row, which already carries the date/text styles
#    this new row needs) and then overwrite the values.
# -----------------------------------------------------------------
$ws.Range("A24:F24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Rows(25).RowHeight = 13.05

$ws.Range("A25").Value = "SCHMITT MUSIC CTR"
$ws.Range("B25").Value = "Monroe, Michael D"
$ws.Range("C25").Value = "003"
$ws.Range("D25").Value = 45954
$ws.Range("E25").Value = "0005169"
